# Updated symbol list refresh (prices / 1h volume %) for cryptos.xlsx.
# Values are set with a leading apostrophe so Excel stores them as literal
# text (matching the original inlineStr cells) instead of re-parsing
# numeric-looking strings ("260.86", "-0.32%", ...) into numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'260.86"
$ws.Range("E2").Value = "'1.83%"
$ws.Range("D3").Value = "'27.40"
$ws.Range("E3").Value = "'2.21%"
$ws.Range("D4").Value = "'4.686"
$ws.Range("E4").Value = "'-0.32%"
$ws.Range("E5").Value = "'2.53%"
$ws.Range("D6").Value = "'6.662"
$ws.Range("E6").Value = "'0.75%"
$ws.Range("D7").Value = "'0.8489"
$ws.Range("E7").Value = "'-0.27%"
$ws.Range("D8").Value = "'0.9285"
$ws.Range("E8").Value = "'1.82%"
$ws.Range("D9").Value = "'0.1404"
$ws.Range("E9").Value = "'1.64%"
$ws.Range("D10").Value = "'0.04881"
$ws.Range("E10").Value = "'9.15%"
$ws.Range("E11").Value = "'1.40%"
$ws.Range("D12").Value = "'0.03078"
$ws.Range("E12").Value = "'0.59%"
$ws.Range("D13").Value = "'0.09061"
$ws.Range("E13").Value = "'-0.29%"
$ws.Range("D14").Value = "'0.001542"
$ws.Range("E14").Value = "'1.32%"
$ws.Range("D15").Value = "'0.0006082"
$ws.Range("E15").Value = "'0.81%"
$ws.Range("D16").Value = "'0.006120"
$ws.Range("E16").Value = "'1.48%"
$ws.Range("E17").Value = "'-0.61%"
$ws.Range("E18").Value = "'-0.38%"
$ws.Range("E19").Value = "'-0.62%"
$ws.Range("E20").Value = "'2.67%"
$ws.Range("E21").Value = "'0.82%"
$ws.Range("D22").Value = "'4.085"
$ws.Range("E22").Value = "'5.42%"
$ws.Range("D23").Value = "'0.04238"
$ws.Range("E23").Value = "'-0.46%"
$ws.Range("D24").Value = "'0.001223"
$ws.Range("E24").Value = "'0.53%"
$ws.Range("E27").Value = "'3.35%"
$ws.Range("D40").Value = "'0.03858"
$ws.Range("E40").Value = "'2.26%"
$ws.Range("D41").Value = "'0.1114"
$ws.Range("E41").Value = "'1.69%"
$ws.Range("D42").Value = "'0.004073"
$ws.Range("E42").Value = "'-34.26%"
$ws.Range("E43").Value = "'15.23%"
$ws.Range("D44").Value = "'0.002217"
$ws.Range("E44").Value = "'0.69%"
$ws.Range("D45").Value = "'0.00005140"
$ws.Range("E45").Value = "'-3.28%"
$ws.Range("E46").Value = "'-0.04%"
$ws.Range("D47").Value = "'0.1371"
$ws.Range("E47").Value = "'-39.28%"
$ws.Range("E48").Value = "'23.71%"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'-0.04%"
$ws.Range("E50").Value = "'-0.04%"